# Update the division-problem answers in the practice table.
#
# The table has 20 rows; the 5 "content" rows (1, 5, 9, 13, 17) each
# hold 5 worked problems ("a÷b=q, r") in their 5 cells, the other rows
# being blank spacer rows. We address each cell positionally via
# Table.Cell(row, col) and overwrite its Range.Text directly.
#
# Note: several old/new values collide with each other across
# different cells (e.g. "50÷6=8, 2" is both the old value of Row9/Col1
# and the new value written into Row17/Col2; "95÷8=11, 7" is both the
# old value of Row9/Col2 and the new value written into Row13/Col5).
# Find/Execute in this runtime always searches from the start of the
# document regardless of which Range it is invoked on, so a text-based
# Find&Replace would mis-target those cells. Assigning Range.Text
# directly (scoped to each cell) avoids that ambiguity entirely while
# still preserving the existing run formatting (font/size) because the
# run's rPr stays attached to the (now retexted) run.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="62÷2=31, 0"},
    @{Row=1;  Col=2; New="75÷6=12, 3"},
    @{Row=1;  Col=3; New="16÷8=2, 0"},
    @{Row=1;  Col=4; New="24÷2=12, 0"},
    @{Row=1;  Col=5; New="64÷5=12, 4"},

    @{Row=5;  Col=1; New="83÷7=11, 6"},
    @{Row=5;  Col=2; New="56÷8=7, 0"},
    @{Row=5;  Col=3; New="23÷6=3, 5"},
    @{Row=5;  Col=4; New="18÷7=2, 4"},
    @{Row=5;  Col=5; New="71÷4=17, 3"},

    @{Row=9;  Col=1; New="50÷6=8, 2"},
    @{Row=9;  Col=2; New="85÷7=12, 1"},
    @{Row=9;  Col=3; New="44÷7=6, 2"},
    @{Row=9;  Col=4; New="88÷4=22, 0"},
    @{Row=9;  Col=5; New="51÷5=10, 1"},

    @{Row=13; Col=1; New="77÷7=11, 0"},
    @{Row=13; Col=2; New="98÷7=14, 0"},
    @{Row=13; Col=3; New="99÷7=14, 1"},
    @{Row=13; Col=4; New="28÷2=14, 0"},
    @{Row=13; Col=5; New="95÷8=11, 7"},

    @{Row=17; Col=1; New="86÷6=14, 2"},
    @{Row=17; Col=2; New="13÷6=2, 1"},
    @{Row=17; Col=3; New="97÷4=24, 1"},
    @{Row=17; Col=4; New="84÷2=42, 0"},
    @{Row=17; Col=5; New="93÷6=15, 3"}
)

foreach ($u in $updates) {
    $rng = $t.Cell($u.Row, $u.Col).Range
    $rng.MoveEnd(12, -1)   # wdCharacter: trim the trailing cell-end mark
    $rng.Text = $u.New
}
